$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $rng.Find.Replacement.ClearFormatting()
    $result = $rng.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $result) {
        Write-Host "WARNING: replace failed for: $old"
    }
}

function Delete-Text($old) {
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $found = $rng.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host "WARNING: delete-find failed for: $old"
    } else {
        $rng.Text = ""
    }
}

# Title
Replace-Text "The Enduring Enigma of Consciousness" "Biology: Exploring the Wonders of Life"

# Author name (merges 3 runs "Dr" + "." + " Anya Gupta" into one run)
Replace-Text "Dr. Anya Gupta" "Helen Robertson"

# Email user part
Replace-Text "agupta@neuro" "helenrobertson@edunet"

# Body paragraph sentences
Replace-Text "The intricate workings of consciousness have long captivated the curious minds of philosophers, scientists, and artists alike" "Our world teems with an astounding array of life, from the tiniest microbes to the majestic blue whales that grace our oceans"

Replace-Text " Its elusive nature has given rise to a plethora of theories and hypotheses, encompassing diverse fields such as neuroscience, psychology, and philosophy" " This captivating tapestry of living organisms, known as biodiversity, holds immense significance for our survival and well-being"

Replace-Text " Despite the advancements in modern science, consciousness remains an enigmatic phenomenon, evoking profound questions about our existence, perception, and relationship with the universe" " Biology, the study of life, unveils the intricate workings of these organisms, delving into their structure, function, growth, and evolution"

Replace-Text " Unraveling the complexities of consciousness presents an intellectual frontier that promises transformative insights into the essence of being" " This exploration guides us toward comprehending the beauty and complexity of the natural world"

Replace-Text "As we delve into the depths of consciousness, we encounter a myriad of perplexing questions" "Biology is a captivating subject that probes the mechanisms underlying the diversity of life"

Replace-Text " What are the physiological and neural correlates of consciousness? How do subjective experiences arise from electrochemical processes in the brain? What is the relationship between consciousness and the physical world we perceive? These questions challenge the boundaries of our understanding and invite us to explore the fundamental nature of reality" " This includes understanding the structure and function of cells, the basic unit of life, and how they work together to form tissues, organs, and organ systems in complex organisms"

Replace-Text " The journey to understanding consciousness promises to shed light on the interconnectedness of mind, body, and the cosmos, offering a glimpse into the deepest mysteries of human existence" " Biology also examines how organisms interact with their environment, revealing their fascinating adaptations and ecological relationships"

Replace-Text "Furthermore, the exploration of consciousness has profound implications for our understanding of free will, moral responsibility, and the nature of qualia" "From the grandeur of a towering forest ecosystem to the microscopic world of bacteria, biology unravels the interconnectedness of all living things"

Replace-Text " By unraveling the mechanisms underlying conscious experience, we may gain insights into the relationship between the conscious and unconscious mind, the genesis of creativity and inspiration, and the essence of self-awareness" " It elucidates the fundamental principles governing reproduction, genetics, and evolution, providing insights into the origin and diversification of life on Earth"

Replace-Text " The quest to comprehend consciousness is a testament to our insatiable curiosity, our desire to understand the universe and our place within it" " Understanding biology equips us with the knowledge to appreciate the fragility of our planet and the importance of preserving its biodiversity"

# Remove the trailing sentence (and its leading period) that no longer exists in the new version
Delete-Text " It is a journey fraught with challenges, but the potential rewards are immense, promising transformative insights into the very essence of existence."

# Summary paragraph sentences
Replace-Text "The study of consciousness, an enduring enigma that has challenged scholars for centuries, presents a captivating frontier of intellectual exploration" "Biology, the study of life, unveils the intricate workings of living organisms, from the microscopic to the majestic"

Replace-Text " From the intricacies of neural processes to the nature of subjective experience, the quest to understand consciousness promises transformative insights into our existence, perception, and relationship with the universe" " It delves into their structure, function, growth, and evolution, revealing the captivating tapestry of biodiversity that sustains our world"

Replace-Text " Its ramifications extend to fundamental questions of free will, moral responsibility, and the essence of qualia" " Through its examination of cells, organisms, and their interactions with the environment, biology equips us with an understanding of the fundamental principles of life, guiding us toward appreciating the interconnectedness of all living things and the importance of preserving our planet's biodiversity"

# Remove the trailing sentences (including the lastRenderedPageBreak run) that no longer exist
Delete-Text " As we delve deeper into the mysteries of the conscious mind, we may unlock the secrets of creativity, inspiration, and self-awareness, gaining a profound understanding of the universe and our place within it."

# Add a new empty paragraph at the end of the document
$d.Paragraphs.Add() | Out-Null

# Apply the "Times New Roman" font to all non-empty paragraph content (excluding the
# paragraph mark itself, to avoid introducing new w:pPr/w:rPr elements)
foreach ($p in $d.Paragraphs) {
    $r = $p.Range
    if (($r.End - 1) -gt $r.Start) {
        $r2 = $d.Range($r.Start, $r.End - 1)
        $r2.Font.Name = "Times New Roman"
    }
}

Write-Host "Edit complete"
